# Daily attendance processing - 2025-11-20 18:55:34
#
# Column G ("Recorded By") holds a comma-separated list of the users who
# touched a given attendance record. This pass reverses the ordering of
# that list for every row that has more than one entry (single-name cells
# are left untouched, since reversing a 1-element list is a no-op anyway).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$usedRange = $ws.UsedRange
$firstRow = $usedRange.Row
$lastRow = $firstRow + $usedRange.Rows.Count - 1

$recordedByCol = 7   # column G

for ($r = $firstRow; $r -le $lastRow; $r++) {
    $cell = $ws.Cells.Item($r, $recordedByCol)
    $val = $cell.Value2

    if ($val -ne $null -and $val -like "*,*") {
        $parts = $val -split ", "
        $count = $parts.Count

        $reversed = @()
        for ($i = $count - 1; $i -ge 0; $i--) {
            $reversed += $parts[$i]
        }

        $cell.Value = $reversed -join ", "
    }
}
